# Apply the edits described by the commit: update the Execute flags on the
# Buses sheet (swap Y/N for rows 2 and 3), and switch the active sheet /
# selection from "Buses" (D24) to "TestData" (C11), leaving a lingering
# selection of D8 on the Buses sheet.

$wb = $excel.ActiveWorkbook

$wsBuses = $wb.Worksheets.Item("Buses")
$wsTestData = $wb.Worksheets.Item("TestData")

# Swap the Execute values for the two test rows.
$wsBuses.Range("B2").Value = "N"
$wsBuses.Range("B3").Value = "Y"

# Leave a selection behind on the Buses sheet (no longer the active tab).
$wsBuses.Range("D8").Select()

# Make TestData the active/selected sheet with its own selection.
$wsTestData.Activate()
$wsTestData.Range("C11").Select()
